$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.432.61"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.100.87"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.81"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5226"
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4575"
$ws.Range("E8").Value = "  +3.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "56.65"
$ws.Range("E9").Value = "  +18.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08924"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.20"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.099.99"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.808"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.053"
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.10"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001147"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06635"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.20"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.296"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.492.23"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.35"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.345.54"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.20"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.38"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.515"
$ws.Range("E29").Value = "  -5.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.26"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.208"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1068"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.381"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.31"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.917"
$ws.Range("E37").Value = "  +6.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02576"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06858"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2324"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.67"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6875"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.246"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.329"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.03"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6383"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.659"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.250"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "83.20"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000339"
$ws.Range("E50").Value = "  +10.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.196"
$ws.Range("E51").Value = "  -1.62%  "
